$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5: "Perkins Wall" -> "Perkins Ceiling" (keeps Neighbor=Y in D5)
$ws.Range("B5").Value = "Perkins Ceiling"

# Row 6: was "Perkins Kitchen Tiles" (Passable=Y in C6) -> now "Perkins Brick Wall" (Neighbor=Y in D6)
$ws.Range("B6").Value = "Perkins Brick Wall"
$ws.Range("C6").ClearContents()
$ws.Range("D6").Value = "Y"

# Row 7: was "Perkins Kitchen Tiles" w/ note in F7 -> now "Perkins Tile Wall", no note
$ws.Range("B7").Value = "Perkins Tile Wall"
$ws.Range("F7").ClearContents()

# Row 8: was "Grass" (Passable=Y) -> now "Perkins Kitchen Tiles" (Passable=Y, unchanged col)
$ws.Range("B8").Value = "Perkins Kitchen Tiles"

# Row 9: was "Impassable Grass" w/ F9 note -> now "Perkins Kitchen Tiles" w/ different F9 note
$ws.Range("B9").Value = "Perkins Kitchen Tiles"
$ws.Range("F9").Value = "for decorations (oven, counter, sink)"

# Row 10: was "Sidewalk" (Passable=Y, Neighbor=Y, F note) -> now "Grass" (Passable=Y only)
$ws.Range("B10").Value = "Grass"
$ws.Range("D10").ClearContents()
$ws.Range("F10").ClearContents()

# Row 11: was "Road" (Passable=Y, Neighbor=Y, F note) -> now "Impassable Grass" (F note only)
$ws.Range("B11").Value = "Impassable Grass"
$ws.Range("C11").ClearContents()
$ws.Range("D11").ClearContents()
$ws.Range("F11").Value = "for use of décor or complex graphics"

# Row 12: was "Dirt" (Passable=Y) -> now "Sidewalk" (Passable=Y, Neighbor=Y, F note)
$ws.Range("B12").Value = "Sidewalk"
$ws.Range("D12").Value = "Y"
$ws.Range("F12").Value = "sides are grass"

# Row 13: was empty (only formula in A13) -> now "Road" (Passable=Y, Neighbor=Y, F note)
$ws.Range("B13").Value = "Road"
$ws.Range("C13").Value = "Y"
$ws.Range("D13").Value = "Y"
$ws.Range("F13").Value = "sides are grass"

# Row 14: was empty (only formula in A14) -> now "Dirt" (Passable=Y)
$ws.Range("B14").Value = "Dirt"
$ws.Range("C14").Value = "Y"

# Update the active selection to match the new cursor position recorded in the file
$ws.Activate()
$ws.Range("E6").Select()
